$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update DESCRIPTION column (B) text to new "Verify that ..." wording ---
$ws.Cells.Item(2,2).Value = "Verify that to test the authorize API for Redirection to TR ID login page"
$ws.Cells.Item(3,2).Value = "Verify that to test the authorize API for Redirection to Facebook login page"
$ws.Cells.Item(4,2).Value = "Verify that to test the authorize API for Redirection to Linked-In login page"
$ws.Cells.Item(5,2).Value = "Verify that user is able to evict a user by passing truid"
$ws.Cells.Item(6,2).Value = "Verify that user is able to activate evicted user by passing truid"
$ws.Cells.Item(7,2).Value = "Verify that get evicted user info or bucket information"
$ws.Cells.Item(8,2).Value = "Verify that user able to get User Mail id by passing truid"

# --- Clear out the STATUS column (L) results for data rows ---
$ws.Cells.Item(2,12).ClearContents()
$ws.Cells.Item(3,12).ClearContents()
$ws.Cells.Item(4,12).ClearContents()
$ws.Cells.Item(5,12).ClearContents()
$ws.Cells.Item(6,12).ClearContents()
$ws.Cells.Item(7,12).ClearContents()
$ws.Cells.Item(8,12).ClearContents()

# --- Clear stray empty placeholder cells that Excel drops on save ---
$ws.Cells.Item(5,7).ClearContents()
$ws.Cells.Item(5,9).ClearContents()
$ws.Cells.Item(5,11).ClearContents()
$ws.Cells.Item(6,7).ClearContents()
$ws.Cells.Item(6,9).ClearContents()
$ws.Cells.Item(6,11).ClearContents()
$ws.Cells.Item(7,6).ClearContents()
$ws.Cells.Item(7,7).ClearContents()
$ws.Cells.Item(7,9).ClearContents()
$ws.Cells.Item(7,11).ClearContents()
$ws.Cells.Item(8,6).ClearContents()
$ws.Cells.Item(8,9).ClearContents()
$ws.Cells.Item(8,11).ClearContents()

# --- Adjust row heights to match updated wrapped text ---
$ws.Rows.Item(3).RowHeight = 45
$ws.Rows.Item(4).RowHeight = 45
$ws.Rows.Item(5).RowHeight = 30
$ws.Rows.Item(6).RowHeight = 30

# --- Update the selected cell / view state ---
$ws.Range("B4").Select()
